$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: new status entry. Column A holds a date string that must be
# stored as literal text (matching the existing A2:A10 date-as-text
# cells) rather than being auto-converted to a date serial number, so
# enter it with a leading apostrophe (Excel's "force text" prefix) and
# then strip the resulting formatting flag back off the cell.
$ws.Range("A11").Value = "'1/24/2010"
$ws.Range("A11").ClearFormats()

$ws.Range("B11").Value = 1.75
$ws.Range("C11").Value = "CFP Update/Edit"

$null = $ws.Range("A12").Select()
